$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = -0.3261774413778852
$ws.Range("J4").Value = 0.4961668977682696
$ws.Range("K4").Value = 0.7648173613827539
$ws.Range("L4").Value = 3.206213021077398
